$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = 44491
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101007
$ws.Cells.Item($row, 10).Value = "Kiwi"
$ws.Cells.Item($row, 11).Value = "Hayward"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 300
$ws.Cells.Item($row, 14).Value = 14000
$ws.Cells.Item($row, 15).Value = 15000
$ws.Cells.Item($row, 16).Value = 14500
$ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 1450
$ws.Cells.Item($row, 20).Value = 10
